$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set column C (Fitness), rows 2 through 252, to the new constant value 7573
$ws.Range("C2:C252").Value = 7573
